# This edit reshuffles the per-observation data across a number of rows in the
# "Artfynd" sheet: the content of several whole rows is rotated/exchanged
# among each other (row numbers, formatting and the header stay untouched -
# only the observation data that lives in each row moves to a different row).
#
# To move a row's data reliably (including clearing cells that should become
# empty, e.g. columns K/L/M/N/AC for some species) we snapshot each involved
# row's values with Range.Value2 and write them back with Range.Value.
#
# Columns Y and AA hold the textual date "2026-01-31" for every single row.
# Excel auto-converts such a string into a real date value as soon as it is
# written back through Value/Value2, which would needlessly change the cell
# type/format even though the visible content never actually changes as part
# of this edit. We therefore simply skip columns Y and AA (splitting each row
# into the A:X, Z and AB:AY pieces) so those two cells are left completely
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowParts([int]$r) {
    $p1 = $ws.Range("A$r" + ":X$r").Value2
    $p2 = $ws.Range("Z$r").Value2
    $p3 = $ws.Range("AB$r" + ":AY$r").Value2
    return @($p1, $p2, $p3)
}

function Set-RowParts([int]$r, $parts) {
    $ws.Range("A$r" + ":X$r").Value = $parts[0]
    $ws.Range("Z$r").Value = $parts[1]
    $ws.Range("AB$r" + ":AY$r").Value = $parts[2]
}

# Maps: destination row number -> row number whose CURRENT data should end up there.
$moves = @{
    2  = 4
    4  = 5
    5  = 2
    7  = 8
    8  = 7
    11 = 12
    12 = 11
    16 = 19
    17 = 18
    18 = 17
    19 = 16
    20 = 22
    21 = 20
    22 = 21
    25 = 28
    26 = 27
    27 = 26
    28 = 25
    30 = 31
    31 = 33
    32 = 34
    33 = 32
    34 = 30
    37 = 38
    38 = 37
    40 = 41
    41 = 42
    42 = 40
    44 = 45
    45 = 44
    47 = 50
    50 = 47
    52 = 57
    53 = 55
    54 = 56
    55 = 54
    56 = 53
    57 = 52
}

# Snapshot the current content of every row involved in the shuffle before
# overwriting anything, since several rows are both a source and a target.
$snapshot = @{}
foreach ($r in $moves.Keys) {
    $snapshot[$r] = Get-RowParts $r
}

foreach ($destRow in $moves.Keys) {
    $srcRow = $moves[$destRow]
    Set-RowParts $destRow $snapshot[$srcRow]
}
